# Add a "Spain" row and a "Total" row to the foreign-flag table, right
# after the existing "Malaysia" row (which was previously the last row).
#
# wdBorder* indices used on Cell.Borders.Item(...):
#   -1 = top, -2 = left, -3 = bottom, -4 = right
# LineWidth maps to half of the OOXML w:sz (w:sz = LineWidth * 2), e.g.
# LineWidth=4 -> w:sz="8", LineWidth=8 -> w:sz="16".
# Row.Height is in points; OOXML w:trHeight (twips) = Height * 20.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$borderColor = 6710886   # 666666
$noColor = 0

function Set-CellBorder($cell, $side, $width, $color, $style) {
    $b = $cell.Borders.Item($side)
    $b.LineWidth = $width
    $b.Color = $color
    $b.LineStyle = $style
}

# --- The current last row ("Malaysia") is no longer the last row, so its
#     heavy bottom rule goes away (same as every other interior body row). ---
$lastRow = $t.Rows.Item($t.Rows.Count)
for ($i = 1; $i -le $lastRow.Cells.Count; $i++) {
    Set-CellBorder $lastRow.Cells.Item($i) -3 0 0 0
}

# --- New row: Spain ---
$spainRow = $t.Rows.Add()
$spainRow.HeightRule = 0
$spainRow.Height = 30.75

$spainValues = @("Spain", "", "0", "1")
for ($i = 1; $i -le $spainRow.Cells.Count; $i++) {
    $cell = $spainRow.Cells.Item($i)
    $cell.Range.Text = $spainValues[$i - 1]
    Set-CellBorder $cell -1 0 0 0
    Set-CellBorder $cell -3 4 $borderColor 1
    Set-CellBorder $cell -2 0 0 0
    Set-CellBorder $cell -4 0 0 0
}

# --- New row: Total ---
$totalRow = $t.Rows.Add()
$totalRow.HeightRule = 0
$totalRow.Height = 30

$totalValues = @("Total", "1,753,792", "4,218", "121")
for ($i = 1; $i -le $totalRow.Cells.Count; $i++) {
    $cell = $totalRow.Cells.Item($i)
    $cell.Range.Text = $totalValues[$i - 1]
    Set-CellBorder $cell -1 4 $borderColor 1
    Set-CellBorder $cell -3 8 $borderColor 1
    Set-CellBorder $cell -2 0 0 0
    Set-CellBorder $cell -4 0 0 0
}

Write-Output ("Final row count=" + $t.Rows.Count)
